$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row 7 (everything currently on/after row 7 shifts down
# by one) and populate it with the new "Gallstone (Large) + Adenomyoma"
# entry under "Gallbladder and biliary tract".
# ------------------------------------------------------------------
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value2 = "Gallbladder and biliary tract"
$ws.Range("B7").Value2 = "Gallstone (Large) + Adenomyoma"
$ws.Range("C7").Value2 = "Clip 1 B-mode"
$ws.Range("D7").Value2 = "https://youtu.be/mnDuOgdSpLA"

# Match the existing "YouTube Link" column formatting (hyperlink style)
# used by every other row in column D.
$ws.Range("D7").Style = $ws.Range("D8").Style

# ------------------------------------------------------------------
# The row insert shifts cell contents down, but this engine does not
# automatically re-point the worksheet's Hyperlinks collection at the
# cells' new addresses. Rebuild the hyperlink list from scratch so
# every link lands on the right (now shifted) row, plus add the brand
# new link for the "Umbilical Vein" row that didn't have one before.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D4"),  "https://youtu.be/zxTC0YBY2RY")
$ws.Hyperlinks.Add($ws.Range("D29"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/91M82AIMyu0")
$ws.Hyperlinks.Add($ws.Range("D34"), "https://youtu.be/qushjTAy6XQ")
$ws.Hyperlinks.Add($ws.Range("D31"), "https://youtu.be/pc-vbxSRTbs")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://youtu.be/DjI1kEnzfSQ")
$ws.Hyperlinks.Add($ws.Range("D30"), "https://youtu.be/JvwODCASLYQ")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/U3ydTsRwxok")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://youtu.be/15o_Km86IzM")
$ws.Hyperlinks.Add($ws.Range("D35"), "https://youtu.be/_FckFwJwynI")
$ws.Hyperlinks.Add($ws.Range("D32"), "https://youtu.be/Axbee4vjNtU")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://youtu.be/RhSUFLTmTl4")
$ws.Hyperlinks.Add($ws.Range("D8"),  "https://youtu.be/2kRZcpi70Aw")
$ws.Hyperlinks.Add($ws.Range("D36"), "https://youtu.be/z_oaRVxRz5s")
$ws.Hyperlinks.Add($ws.Range("D5"),  "https://youtu.be/K2Wbg7BgXy4")
$ws.Hyperlinks.Add($ws.Range("D3"),  "https://youtu.be/ZXwd0gwHEkQ")
$ws.Hyperlinks.Add($ws.Range("D33"), "https://youtu.be/VJdnjrAAO-4")
$ws.Hyperlinks.Add($ws.Range("D2"),  "https://youtu.be/kdZO1IPuOIw")
$ws.Hyperlinks.Add($ws.Range("D37"), "https://youtu.be/S45odD2wQOQ")
$ws.Hyperlinks.Add($ws.Range("D27"), "https://youtu.be/ytNgK7wuL_M")

# Restore the new link cell's own hyperlink style again (Hyperlinks.Add
# stamps a fresh style xf; keep it aligned with the rest of column D).
$ws.Range("D7").Style = $ws.Range("D8").Style

# ------------------------------------------------------------------
# Misc view state that Excel also recorded for this edit.
# ------------------------------------------------------------------
$ws.Range("D7").Select()
